{"js": "// Left-align the first column of every table: every first-column\n// paragraph that is currently right-aligned becomes left-aligned\n// (the header row's first cell is already left-aligned, so it is\n// naturally skipped). Other columns are left untouched.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (const table of tables.items) {\n  table.rows.load(\"items\");\n}\nawait context.sync();\n\nfor (const table of tables.items) {\n  for (const row of table.rows.items) {\n    row.cells.load(\"items\");\n  }\n}\nawait context.sync();\n\nconst firstCellParagraphs = [];\n\nfor (const table of tables.items) {\n  for (const row of table.rows.items) {\n    const firstCell = row.cells.items[0];\n    firstCell.body.paragraphs.load(\"items\");\n    firstCellParagraphs.push(firstCell.body.paragraphs);\n  }\n}\nawait context.sync();\n\nfor (const paragraphs of firstCellParagraphs) {\n  for (const paragraph of paragraphs.items) {\n    paragraph.load(\"alignment\");\n  }\n}\nawait context.sync();\n\nfor (const paragraphs of firstCellParagraphs) {\n  for (const paragraph of paragraphs.items) {\n    if (paragraph.alignment === Word.Alignment.right) {\n      paragraph.alignment = Word.Alignment.left;\n    }\n  }\n}\nawait context.sync();\n", "ps1": "# Left-align the first column of every table: every first-column\n# paragraph that is currently right-aligned becomes left-aligned\n# (the header row's first cell is already left-aligned, so it is\n# naturally skipped). Other columns are left untouched.\n\n$d = $word.ActiveDocument\n\nforeach ($table in $d.Tables) {\n    foreach ($row in $table.Rows) {\n        $cell = $row.Cells.Item(1)\n        if ($cell.Range.ParagraphFormat.Alignment -eq 2) {\n            $cell.Range.ParagraphFormat.Alignment = 0\n        }\n    }\n}\n"}
